# Rapport du 08 Octobre 2025
# Add 8 new delivery rows (week of 07 Oct 2025) to the bottom of the
# "Semaine_1" table, mirroring how a user would type new rows directly
# below an existing Excel Table (auto-growing it from A1:P27 to A1:P35).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

$firstNewRow = 28
$lastNewRow  = 35

# 1) Grow the table range so it covers the new rows. Using Resize (rather
#    than ListRows.Add in a loop) leaves the existing rows completely
#    untouched - only the table/autofilter ref and sheet dimension change.
$lo.Resize($ws.Range("A1:P" + $lastNewRow))

# 2) The row that used to be the table's last row (27) carries a slightly
#    different ("last row") cell format than the rest of the body rows.
#    Once it's no longer the last row it should pick up the regular body
#    format instead - copy formats (only) from the row above (26) across
#    every row from the old last row down through the new last row so
#    they all share the same regular body formatting/font.
$ws.Range("A26:P26").Copy() | Out-Null
$ws.Range(("A{0}:P{1}" -f ($firstNewRow - 1), $lastNewRow)).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# 3) New row data (Prenom_Nom_RZ / zone / secteur / Nom_du_magasin /
#    Telephone_Client / Type / Point_de_Vente / Operation / Commentaire /
#    Produit / Quantites / Prix_Unitaire / Prix Total). Date + calculated
#    Semaine/Mois columns are set separately below.
$rows = @(
    @{ B="Seynabou SOW"; C="CASTOR"; D="Castor";          E="Ndioguou";        F=776634479; G="Grossiste"; H="Client Partenaire";     I="Livraison"; J="Merci beaucoup";       K="Café pot Refraish 50g";   L=1;  M=10250; N=10250  },
    @{ B="Seynabou SOW"; C="CASTOR"; D="Zone de captage";  E="Alpha";           F=704738232; G="Demi-Gros"; H="Client non Partenaire"; I="Livraison"; J="Je vais essayer avec "; K="Kamlac évaporé 48x160g";  L=1;  M=11500; N=11500  },
    @{ B="Seynabou SOW"; C="CASTOR"; D="Castor";          E="Mor Diop";        F=777262311; G="Grossiste"; H="Client Partenaire";     I="Livraison"; J="Je vais essayer ";      K="Kamlac évaporé 48x160g";  L=2;  M=11500; N=23000  },
    @{ B="Seynabou SOW"; C="CASTOR"; D="Castor";          E="Assane Wade";     F=775884054; G="Demi-Gros"; H="Client non Partenaire"; I="Livraison"; J="Je vais essayer avec "; K="Kamlac évaporé 48x160g";  L=2;  M=11500; N=23000  },
    @{ B="Seynabou SOW"; C="CASTOR"; D="Liberté 1 à 6";    E="Moussa";          F=771837885; G="Demi-Gros"; H="Client Partenaire";     I="Livraison"; J="Je vais essayer avec "; K="Kamlac évaporé 48x160g";  L=2;  M=11500; N=23000  },
    @{ B="Seynabou SOW"; C="CASTOR"; D="Liberté 1 à 6";    E="Omar";            F=773170826; G="Demi-Gros"; H="Client non Partenaire"; I="Livraison"; J="Je vais essayer ";      K="Kamlac évaporé 48x160g";  L=1;  M=11500; N=11500  },
    @{ B="Seynabou SOW"; C="CASTOR"; D="Ngor";             E="Cheikh Boussole"; F=777802399; G="Demi-Gros"; H="Client non Partenaire"; I="Livraison"; J="Merci beaucoup ";      K="Kamlac évaporé 48x160g";  L=3;  M=11500; N=34500  },
    @{ B="Seynabou SOW"; C="CASTOR"; D="Zone de captage";  E="El Hadji";        F=773739328; G="Demi-Gros"; H="Client Partenaire";     I="Livraison"; J="Je vais essayer avec "; K="Kamlac évaporé 48x160g";  L=12; M=11500; N=138000 }
)

$r = $firstNewRow
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value  = 45937          # Date (07/10/2025)
    $ws.Cells.Item($r, 2).Value  = $row.B          # Prenom_Nom_RZ
    $ws.Cells.Item($r, 3).Value  = $row.C          # zone
    $ws.Cells.Item($r, 4).Value  = $row.D          # secteur
    $ws.Cells.Item($r, 5).Value  = $row.E          # Nom_du_magasin
    $ws.Cells.Item($r, 6).Value  = $row.F          # Telephone_Client
    $ws.Cells.Item($r, 7).Value  = $row.G          # Type
    $ws.Cells.Item($r, 8).Value  = $row.H          # Point_de_Vente
    $ws.Cells.Item($r, 9).Value  = $row.I          # Operation
    $ws.Cells.Item($r, 10).Value = $row.J          # Commentaire
    $ws.Cells.Item($r, 11).Value = $row.K          # Produit
    $ws.Cells.Item($r, 12).Value = $row.L          # Quantites
    $ws.Cells.Item($r, 13).Value = $row.M          # Prix_Unitaire
    $ws.Cells.Item($r, 14).Value = $row.N          # Prix Total
    $ws.Cells.Item($r, 15).Formula = '="S"&_xlfn.ISOWEEKNUM(Semaine_1[[#This Row],[Date]])'
    $ws.Cells.Item($r, 16).Formula = '=TEXT(Semaine_1[[#This Row],[Date]],"MMMM")'
    $r++
}

# 4) Match the view state left behind by the edit: scrolled back so the
#    top-left cell is the sheet default again, with E11 selected.
$ws.Range("E11").Select()

$wb.Save()
